# Duplicate the "20233" sheet into a brand-new trailing sheet named "20234"
# (same data, styles, and merged cells -- mirrors how the workbook's other
# year/quarter tabs were produced).

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("20233")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy the source sheet so it lands right after the current last sheet.
$source.Copy($null, $lastSheet)

# The freshly-copied sheet is now the last tab; give it its proper name.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "20234"
